$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "7.27") need the
# column format forced to Text first, otherwise Excel COM auto-converts the
# assigned string into a numeric value (standard Excel behavior), which would
# change the cell's stored type away from the original text/inlineStr semantics.
$ws.Range("D2").Value = '58.320.32'
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").Value = '3.142.22'
$ws.Range("E3").Value = '  -3.61%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.34'
$ws.Range("E5").Value = '  -4.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.91'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.138.44'
$ws.Range("E8").Value = '  -3.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.444'
$ws.Range("E9").Value = '  -4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.27'
$ws.Range("E10").Value = '  -6.74%  '
$ws.Range("E11").Value = '  -7.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.378'
$ws.Range("E12").Value = '  -5.79%  '
$ws.Range("D13").Value = '3.675.20'
$ws.Range("E13").Value = '  -4.05%  '
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.46'
$ws.Range("E15").Value = '  -3.99%  '
$ws.Range("D16").Value = '3.135.10'
$ws.Range("D17").Value = '58.274.51'
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("E18").Value = '  -5.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.76'
$ws.Range("E19").Value = '  -4.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.03'
$ws.Range("E20").Value = '  -4.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.91'
$ws.Range("E21").Value = '  -6.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '343.85'
$ws.Range("E22").Value = '  -7.22%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.509'
$ws.Range("E24").Value = '  -3.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.69'
$ws.Range("E25").Value = '  -6.81%  '
$ws.Range("D26").Value = '3.263.63'
$ws.Range("E26").Value = '  -4.24%  '
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").Value = '0.0₃0954'
$ws.Range("E28").Value = '  -5.14%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.86'
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.25'
$ws.Range("E32").Value = '  +2.41%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.86'
$ws.Range("E33").Value = '  -7.00%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.92'
$ws.Range("E34").Value = '  -6.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '21.42'
$ws.Range("E35").Value = '  -4.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.85'
$ws.Range("E36").Value = '  -3.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.70'
$ws.Range("E37").Value = '  -5.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.25'
$ws.Range("E38").Value = '  -4.67%  '
$ws.Range("E39").Value = '  -9.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0686'
$ws.Range("E40").Value = '  -4.40%  '
$ws.Range("D41").Value = '3.169.27'
$ws.Range("E41").Value = '  -3.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.47'
$ws.Range("E42").Value = '  -1.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.86'
$ws.Range("E43").Value = '  -7.38%  '
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.692'
$ws.Range("E45").Value = '  -6.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.92'
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.45'
$ws.Range("E48").Value = '  -6.70%  '
$ws.Range("D49").Value = '2.279.97'
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.15'
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.70'
$ws.Range("E51").Value = '  -1.44%  '
